$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{}

$data[2] = @("ECs", "Mrc1", "Ptprc", "ECs", 3, 1, 14.71647433333334, 44.14942300000001, 0.1054330184450109, 0.1054330184450109, 3, 1, 205.313027, 615.9390810000001, 0.435242422384838, 0.435242422384838, 3021.483892144475, 27193.35502930027, 0.04588892234735186, 0.04588892234735186)
$data[3] = @("ECs", "Mrc1", "Ptprc", "FAPs", 3, 1, 14.71647433333334, 44.14942300000001, 0.1054330184450109, 0.1054330184450109, 3, 1, 0.186821, 0.5604629999999999, 0.0003960412341120362, 0.0003960412341120362, 2.749346451427667, 24.744118062849, 0.00004175582274111919, 0.00004175582274111919)
$data[4] = @("ECs", "Mrc1", "Ptprc", "M2", 3, 1, 14.71647433333334, 44.14942300000001, 0.1054330184450109, 0.1054330184450109, 3, 1, 266.1765593333333, 798.529678, 0.5642668278730386, 0.5642668278730386, 3917.180503563978, 35254.6245320758, 0.05949235487104589, 0.05949235487104589)
$data[5] = @("ECs", "Mrc1", "Ptprc", "sCs", 3, 1, 14.71647433333334, 44.14942300000001, 0.1054330184450109, 0.1054330184450109, 2, 0.6666666666666666, 0.044676, 0.134028, 0.00009470850801135487, 0.00009470850801135488, 0.6574732073160001, 5.917258865844001, 0.000009985403872060642, 0.000009985403872060644)
$data[6] = @("FAPs", "Mrc1", "Ptprc", "ECs", 3, 1, 0.07385333333333334, 0.22156, 0.0005291063388682706, 0.0005291063388682705, 3, 1, 205.313027, 615.9390810000001, 0.435242422384838, 0.435242422384838, 15.16305142070667, 136.46746278636, 0.0002302895246281991, 0.000230289524628199)
$data[7] = @("FAPs", "Mrc1", "Ptprc", "FAPs", 3, 1, 0.07385333333333334, 0.22156, 0.0005291063388682706, 0.0005291063388682705, 3, 1, 0.186821, 0.5604629999999999, 0.0003960412341120362, 0.0003960412341120362, 0.01379735358666667, 0.12417618228, 0.0000002095479274218911, 0.000000209547927421891)
$data[8] = @("FAPs", "Mrc1", "Ptprc", "M2", 3, 1, 0.07385333333333334, 0.22156, 0.0005291063388682706, 0.0005291063388682705, 3, 1, 266.1765593333333, 798.529678, 0.5642668278730386, 0.5642668278730386, 19.65802616196445, 176.92223545768, 0.0002985571554407161, 0.000298557155440716)
$data[9] = @("FAPs", "Mrc1", "Ptprc", "sCs", 3, 1, 0.07385333333333334, 0.22156, 0.0005291063388682706, 0.0005291063388682705, 2, 0.6666666666666666, 0.044676, 0.134028, 0.00009470850801135487, 0.00009470850801135488, 0.003299471520000001, 0.02969524368, 0.00000005011087193356424, 0.00000005011087193356424)
$data[10] = @("M2", "Mrc1", "Ptprc", "ECs", 3, 1, 124.584409, 373.753227, 0.8925582306290469, 0.8925582306290469, 3, 1, 205.313027, 615.9390810000001, 0.435242422384838, 0.435242422384838, 25578.80212879605, 230209.2191591644, 0.3884792064185113, 0.3884792064185113)
$data[11] = @("M2", "Mrc1", "Ptprc", "FAPs", 3, 1, 124.584409, 373.753227, 0.8925582306290469, 0.8925582306290469, 3, 1, 0.186821, 0.5604629999999999, 0.0003960412341120362, 0.0003960412341120362, 23.274983873789, 209.474854864101, 0.0003534898631751831, 0.0003534898631751831)
$data[12] = @("M2", "Mrc1", "Ptprc", "M2", 3, 1, 124.584409, 373.753227, 0.8925582306290469, 0.8925582306290469, 3, 1, 266.1765593333333, 798.529678, 0.5642668278730386, 0.5642668278730386, 33161.44933419677, 298453.0440077709, 0.5036410014890244, 0.5036410014890244)
$data[13] = @("M2", "Mrc1", "Ptprc", "sCs", 3, 1, 124.584409, 373.753227, 0.8925582306290469, 0.8925582306290469, 2, 0.6666666666666666, 0.044676, 0.134028, 0.00009470850801135487, 0.00009470850801135488, 5.565933056484001, 50.09339750835601, 0.00008453285833613182, 0.00008453285833613184)
$data[14] = @("sCs", "Mrc1", "Ptprc", "ECs", 3, 1, 0.2065306666666666, 0.6195919999999999, 0.001479644587073792, 0.001479644587073792, 3, 1, 205.313027, 615.9390810000001, 0.435242422384838, 0.435242422384838, 42.40343634166133, 381.630927074952, 0.0006440040943466107, 0.0006440040943466107)
$data[15] = @("sCs", "Mrc1", "Ptprc", "FAPs", 3, 1, 0.2065306666666666, 0.6195919999999999, 0.001479644587073792, 0.001479644587073792, 3, 1, 0.186821, 0.5604629999999999, 0.0003960412341120362, 0.0003960412341120362, 0.03858426567733333, 0.3472583910959999, 0.0000005860002683118988, 0.0000005860002683118988)
$data[16] = @("sCs", "Mrc1", "Ptprc", "M2", 3, 1, 0.2065306666666666, 0.6195919999999999, 0.001479644587073792, 0.001479644587073792, 3, 1, 266.1765593333333, 798.529678, 0.5642668278730386, 0.5642668278730386, 54.97362225015288, 494.7626002513759, 0.0008349143575276408, 0.0008349143575276408)
$data[17] = @("sCs", "Mrc1", "Ptprc", "sCs", 3, 1, 0.2065306666666666, 0.6195919999999999, 0.001479644587073792, 0.001479644587073792, 2, 0.6666666666666666, 0.044676, 0.134028, 0.00009470850801135487, 0.00009470850801135488, 0.009226964063999999, 0.08304267657599999, 0.0000001401349312288361, 0.0000001401349312288361)

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
